$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-update job: append the newest Pick 3 draw result as a new row.
# The sheet keeps every column as plain text (dates/numbers included), so
# force text formatting before writing the values - otherwise Excel would
# reinterpret "2025-09-27" as a date and "250927" as a number.
$ws.Range("A11:E11").NumberFormat = "@"

$ws.Range("A11").Value = "2025-09-27"
$ws.Range("B11").Value = "Pick 3"
$ws.Range("C11").Value = "250927"
$ws.Range("D11").Value = "2-0-8"
$ws.Range("E11").Value = "2025-09-27T21:34:59.193+04:00"
